$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.163986
$ws.Range("H2").Value = 57.491958
$ws.Range("I2").Value = 0.1197574615923936
$ws.Range("J2").Value = 0.1197574615923936
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 3056.342193085078
$ws.Range("R2").Value = 27507.0797377657
$ws.Range("S2").Value = 0.03572706489743428
$ws.Range("T2").Value = 0.03572706489743428

# Row 3
$ws.Range("G3").Value = 19.163986
$ws.Range("H3").Value = 57.491958
$ws.Range("I3").Value = 0.1197574615923936
$ws.Range("J3").Value = 0.1197574615923936
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 3306.908620689042
$ws.Range("R3").Value = 29762.17758620138
$ws.Range("S3").Value = 0.03865605728591056
$ws.Range("T3").Value = 0.03865605728591057

# Row 4
$ws.Range("G4").Value = 19.163986
$ws.Range("H4").Value = 57.491958
$ws.Range("I4").Value = 0.1197574615923936
$ws.Range("J4").Value = 0.1197574615923936
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 1425.564949968121
$ws.Range("R4").Value = 12830.08454971309
$ws.Range("S4").Value = 0.01666411948185965
$ws.Range("T4").Value = 0.01666411948185966

# Row 5
$ws.Range("G5").Value = 19.163986
$ws.Range("H5").Value = 57.491958
$ws.Range("I5").Value = 0.1197574615923936
$ws.Range("J5").Value = 0.1197574615923936
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 1119.456851279399
$ws.Range("R5").Value = 10075.11166151459
$ws.Range("S5").Value = 0.01308587358641462
$ws.Range("T5").Value = 0.01308587358641462

# Row 6
$ws.Range("G6").Value = 19.163986
$ws.Range("H6").Value = 57.491958
$ws.Range("I6").Value = 0.1197574615923936
$ws.Range("J6").Value = 0.1197574615923936
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 1336.615506976976
$ws.Range("R6").Value = 12029.53956279278
$ws.Range("S6").Value = 0.01562434634077447
$ws.Range("T6").Value = 0.01562434634077447

# Row 7
$ws.Range("G7").Value = 24.05875033333333
$ws.Range("H7").Value = 72.176251
$ws.Range("I7").Value = 0.150345281456851
$ws.Range("J7").Value = 0.1503452814568511
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 3836.977012854546
$ws.Range("R7").Value = 34532.79311569091
$ws.Range("S7").Value = 0.04485228357556557
$ws.Range("T7").Value = 0.04485228357556557

# Row 8
$ws.Range("G8").Value = 24.05875033333333
$ws.Range("H8").Value = 72.176251
$ws.Range("I8").Value = 0.150345281456851
$ws.Range("J8").Value = 0.1503452814568511
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 4151.541797218249
$ws.Range("R8").Value = 37363.87617496424
$ws.Range("S8").Value = 0.04852938376769598
$ws.Range("T8").Value = 0.04852938376769599

# Row 9
$ws.Range("G9").Value = 24.05875033333333
$ws.Range("H9").Value = 72.176251
$ws.Range("I9").Value = 0.150345281456851
$ws.Range("J9").Value = 0.1503452814568511
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 1789.675238503818
$ws.Range("R9").Value = 16107.07714653437
$ws.Range("S9").Value = 0.02092038108037114
$ws.Range("T9").Value = 0.02092038108037114

# Row 10
$ws.Range("G10").Value = 24.05875033333333
$ws.Range("H10").Value = 72.176251
$ws.Range("I10").Value = 0.150345281456851
$ws.Range("J10").Value = 0.1503452814568511
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 1405.382622063621
$ws.Range("R10").Value = 12648.44359857259
$ws.Range("S10").Value = 0.01642819847129456
$ws.Range("T10").Value = 0.01642819847129457

# Row 11
$ws.Range("G11").Value = 24.05875033333333
$ws.Range("H11").Value = 72.176251
$ws.Range("I11").Value = 0.150345281456851
$ws.Range("J11").Value = 0.1503452814568511
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 1678.006797438739
$ws.Range("R11").Value = 15102.06117694865
$ws.Range("S11").Value = 0.01961503456192376
$ws.Range("T11").Value = 0.01961503456192377

# Row 12
$ws.Range("G12").Value = 61.341815
$ws.Range("H12").Value = 184.025445
$ws.Range("I12").Value = 0.3833304853108436
$ws.Range("J12").Value = 0.3833304853108436
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 9783.015776828428
$ws.Range("R12").Value = 88047.14199145584
$ws.Range("S12").Value = 0.1143584119416184
$ws.Range("T12").Value = 0.1143584119416184

# Row 13
$ws.Range("G13").Value = 61.341815
$ws.Range("H13").Value = 184.025445
$ws.Range("I13").Value = 0.3833304853108436
$ws.Range("J13").Value = 0.3833304853108436
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 10585.05139965205
$ws.Range("R13").Value = 95265.46259686848
$ws.Range("S13").Value = 0.1237337949767719
$ws.Range("T13").Value = 0.1237337949767719

# Row 14
$ws.Range("G14").Value = 61.341815
$ws.Range("H14").Value = 184.025445
$ws.Range("I14").Value = 0.3833304853108436
$ws.Range("J14").Value = 0.3833304853108436
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 4563.076879279118
$ws.Range("R14").Value = 41067.69191351206
$ws.Range("S14").Value = 0.05334001675821149
$ws.Range("T14").Value = 0.05334001675821151

# Row 15
$ws.Range("G15").Value = 61.341815
$ws.Range("H15").Value = 184.025445
$ws.Range("I15").Value = 0.3833304853108436
$ws.Range("J15").Value = 0.3833304853108436
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 3583.258465731682
$ws.Range("R15").Value = 32249.32619158513
$ws.Range("S15").Value = 0.04188644453461987
$ws.Range("T15").Value = 0.04188644453461988

# Row 16
$ws.Range("G16").Value = 61.341815
$ws.Range("H16").Value = 184.025445
$ws.Range("I16").Value = 0.3833304853108436
$ws.Range("J16").Value = 0.3833304853108436
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 4278.35947882204
$ws.Range("R16").Value = 38505.23530939836
$ws.Range("S16").Value = 0.05001181709962187
$ws.Range("T16").Value = 0.05001181709962189

# Row 17
$ws.Range("G17").Value = 7.095824666666666
$ws.Range("H17").Value = 21.287474
$ws.Range("I17").Value = 0.04434244264135302
$ws.Range("J17").Value = 0.04434244264135302
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 1131.667927719588
$ws.Range("R17").Value = 10185.01134947629
$ws.Range("S17").Value = 0.01322861477600823
$ws.Range("T17").Value = 0.01322861477600823

# Row 18
$ws.Range("G18").Value = 7.095824666666666
$ws.Range("H18").Value = 21.287474
$ws.Range("I18").Value = 0.04434244264135302
$ws.Range("J18").Value = 0.04434244264135302
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 1224.444839455526
$ws.Range("R18").Value = 11020.00355509973
$ws.Range("S18").Value = 0.01431312905391623
$ws.Range("T18").Value = 0.01431312905391623

# Row 19
$ws.Range("G19").Value = 7.095824666666666
$ws.Range("H19").Value = 21.287474
$ws.Range("I19").Value = 0.04434244264135302
$ws.Range("J19").Value = 0.04434244264135302
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 527.8421167662731
$ws.Range("R19").Value = 4750.579050896458
$ws.Range("S19").Value = 0.006170202277733884
$ws.Range("T19").Value = 0.006170202277733885

# Row 20
$ws.Range("G20").Value = 7.095824666666666
$ws.Range("H20").Value = 21.287474
$ws.Range("I20").Value = 0.04434244264135302
$ws.Range("J20").Value = 0.04434244264135302
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 414.4998612802868
$ws.Range("R20").Value = 3730.498751522582
$ws.Range("S20").Value = 0.004845289731445359
$ws.Range("T20").Value = 0.00484528973144536

# Row 21
$ws.Range("G21").Value = 7.095824666666666
$ws.Range("H21").Value = 21.287474
$ws.Range("I21").Value = 0.04434244264135302
$ws.Range("J21").Value = 0.04434244264135302
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 494.9069198994613
$ws.Range("R21").Value = 4454.162279095152
$ws.Range("S21").Value = 0.00578520680224931
$ws.Range("T21").Value = 0.005785206802249311

# Row 22
$ws.Range("G22").Value = 48.362939
$ws.Range("H22").Value = 145.088817
$ws.Range("I22").Value = 0.3022243289985588
$ws.Range("J22").Value = 0.3022243289985588
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 7713.097423850124
$ws.Range("R22").Value = 69417.87681465111
$ws.Range("S22").Value = 0.09016213329959934
$ws.Range("T22").Value = 0.09016213329959934

# Row 23
$ws.Range("G23").Value = 48.362939
$ws.Range("H23").Value = 145.088817
$ws.Range("I23").Value = 0.3022243289985588
$ws.Range("J23").Value = 0.3022243289985588
$ws.Range("M23").Value = 172.558497
$ws.Range("N23").Value = 517.675491
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 8345.436064342683
$ws.Range("R23").Value = 75108.92457908414
$ws.Range("S23").Value = 0.09755384607873321
$ws.Range("T23").Value = 0.09755384607873321

# Row 24
$ws.Range("G24").Value = 48.362939
$ws.Range("H24").Value = 145.088817
$ws.Range("I24").Value = 0.3022243289985588
$ws.Range("J24").Value = 0.3022243289985588
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 3597.608071506954
$ws.Range("R24").Value = 32378.47264356259
$ws.Range("S24").Value = 0.04205418402987197
$ws.Range("T24").Value = 0.04205418402987197

# Row 25
$ws.Range("G25").Value = 48.362939
$ws.Range("H25").Value = 145.088817
$ws.Range("I25").Value = 0.3022243289985588
$ws.Range("J25").Value = 0.3022243289985588
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 2825.102429711526
$ws.Range("R25").Value = 25425.92186740373
$ws.Range("S25").Value = 0.03302399125221033
$ws.Range("T25").Value = 0.03302399125221034

# Row 26
$ws.Range("G26").Value = 48.362939
$ws.Range("H26").Value = 145.088817
$ws.Range("I26").Value = 0.3022243289985588
$ws.Range("J26").Value = 0.3022243289985588
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 3373.131989888825
$ws.Range("R26").Value = 30358.18790899942
$ws.Range("S26").Value = 0.03943017433814389
$ws.Range("T26").Value = 0.0394301743381439
